$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(44418, 2, 18, 74.62067821905315),
    @(44419, 1, 17, 70.47508498466131),
    @(44420, 3, 20, 82.91186468783683),
    @(44421, 11, 29, 120.2222037973634),
    @(44422, 3, 30, 124.3677970317553),
    @(44423, 1, 22, 91.20305115662052),
    @(44424, 1, 22, 91.20305115662052),
    @(44425, 4, 24, 99.4942376254042),
    @(44426, 0, 23, 95.34864439101236),
    @(44427, 1, 21, 87.05745792222866),
    @(44428, 4, 14, 58.03830528148578),
    @(44429, 0, 11, 45.60152557831026),
    @(44430, 7, 17, 70.47508498466131),
    @(44431, 0, 16, 66.32949175026947)
)

$startRow = 344
$lastRow = $startRow - 1

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($lastRow, 1).Copy() | Out-Null
    $ws.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}
